# dstu2-r4-table.csv sheet update:
#  - insert a "Version Deprecated" column (D), shifting the old "Version Retired"
#    (old D) to E and "Notes" (old E) to F
#  - append two new tracking columns: "current_fmm" (G) and "proposed_fmm" (H)
#  - add two new profile rows (US Core FamilyMemberHistory Profile,
#    US Core PMO ServiceRequest Profile) ahead of the "**Extensions**" marker row
#  - add two new extension rows (US Core Individual Sex Extension,
#    US Core FamilyMemberHistory Recorder Extension) at the end of the table
#  - rename "US Core Implantable Device Profile" -> "US Core Device Profile"
#    and record the rename in its Notes
#  - record deprecations/supersessions for several extensions (Gender Identity,
#    Authentication Time, Sex Extension) and fmm values throughout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table gains 4 rows overall (77 -> 81): insert 4 blank rows ahead of the
# old row 63 ("**Extensions**") so the full data block below can be written
# straight across A1:H81 without clobbering rows that need to shift down.
$ws.Range("A63:A66").EntireRow.Insert()

$rows = @(
    @('DSTU2 Profile', 'Equivalent US Core Profile', 'Version Added', 'Version Deprecated', 'Version Retired', 'Notes', 'current_fmm', 'proposed_fmm'),
    @('Argonaut AllergyIntolerance Profile', 'US Core AllergyIntolerance Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut CarePlan Profile', 'US Core CarePlan Profile', '3.1.1*', '', '', 'Type changed from CarePlan in DSTU2 to CareTeam in R4', 3, 4),
    @('Argonaut CareTeam Profile', 'US Core CareTeam Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Condition Profile', 'US Core Condition Profile', '3.1.1*', '', '5.0.0', 'In version 5.0.0 split US Core Condition Profilee into US Core Condition Encounter Diagnosis Profile and US Core Condition Problems and Health Concerns Profile', 3, 5),
    @('     -', 'US Core Condition Encounter Diagnosis Profile', '5.0.0', '', '', 'see above', 3, 4),
    @('     -', 'US Core Condition Problems and Health Concerns Profile', '5.0.0', '', '', 'see above', 3, 4),
    @('Argonaut Device Profile', 'US Core Device Profile', '3.1.1*', '', '', 'in version 9.0.0 removed the term "implantable" form the official ( canonical ) URL,  name,  and title.', 3, 4),
    @('Argonaut DiagnosticReport Profile', 'US Core DiagnosticReport Profile for Laboratory Results Reporting', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut DocumentReference Profile', 'US Core DocumentReference Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Goal Profile', 'US Core Goal Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Immunization Profile', 'US Core Immunization Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Medication Profile', 'US Core Medication Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut MedicationOrder Profile', 'US Core MedicationRequest Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut MedicationStatement Profile', '     -', '', '', '', 'US Core MedicationStatement Profile not defined in 3.1.1*', 3, 5),
    @('Argonaut Observation Results Profile', 'US Core Laboratory Result Observation Profile', '3.1.1*', '', '', 'In Version 6.0.0 derived from US Core Observation Clinical Test Result Profile', 3, 5),
    @('Argonaut Patient Profile', 'US Core Patient Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Procedure Profile', 'US Core Procedure Profile', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Smoking Status Observation Profile', 'US Core Smoking Status Observation Profile', '3.1.1*', '', '', '', 3, 4),
    @('Argonaut Vital Signs Observation Profile', 'US Core Vital Signs Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Blood Pressure Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core BMI Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Head Circumference Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Body Height Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Body Weight Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Body Temperature Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Heart Rate Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Respiratory Rate Profile', '4.0.0', '', '', 'Previous to version 4.0.0, referenced the FHIR Vital Signs Profile', 3, 5),
    @('     -', 'US Core Pediatric BMI for Age Observation Profile', '3.1.1*', '', '', '3.1.1* based on the FHIR Vital Signs Profile, 4.0.0+ based on the US Core Vital Signs Profile', 3, 5),
    @('     -', 'US Core Pediatric Head Occipital Frontal Circumference Percentile Profile', '3.1.1*', '', '', '3.1.1* based on the FHIR Vital Signs Profile, 4.0.0+ based on the US Core Vital Signs Profile', 3, 5),
    @('     -', 'US Core Pediatric Weight for Height Observation Profile', '3.1.1*', '', '', '3.1.1* based on the FHIR Vital Signs Profile, 4.0.0+ based on the US Core Vital Signs Profile', 3, 5),
    @('     -', 'US Core Pulse Oximetry Profile', '3.1.1*', '', '', '3.1.1* based on the FHIR Vital Signs Profile, 4.0.0+ based on the US Core Vital Signs Profile', 3, 5),
    @('     -', 'US Core DiagnosticReport Profile for Report and Note Exchange', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Encounter Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Location Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Organization Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Practitioner Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core PractitionerRole Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Provenance Profile', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Observation Sexual Orientation Profile', '5.0.0', '', '', 'After version 7.0.0 this extension is no longer a USCDI requirement.', 3, 3),
    @('     -', 'US Core Observation Social History Profile', '5.0.0', '', '6.0.0', 'Superseded by US Core Simple Observation Profile in version 6.0.0', 3, 'NA'),
    @('     -', 'US Core Observation Survey Profile', '5.0.0', '', '6.0.0', 'Superseded by US Core Observation Screening Assessment Profile in version 6.0.0', 3, 'NA'),
    @('     -', 'US Core Observation SDOH Assessment Profile', '5.0.0', '', '6.0.0', 'Superseded by US Core Observation Screening Assessment Profile in version 6.0.0', 3, 'NA'),
    @('     -', 'US Core Observation Clinical Test Result Profile', '5.0.0', '', '6.0.0', 'Superseded by US Core Observation Clinical Result Profile in version 6.0.0', 3, 'NA'),
    @('     -', 'US Core Observation Imaging Result Profile', '5.0.0', '', '6.0.0', 'Superseded by US Core Observation Clinical Result Profile in version 6.0.0', 3, 'NA'),
    @('     -', 'US Core QuestionnaireResponse Profile', '5.0.0', '', '', '', 3, 3),
    @('     -', 'US Core RelatedPerson Profile', '5.0.0', '', '', '', 3, 4),
    @('     -', 'US Core ServiceRequest Profile', '5.0.0', '', '', '', 3, 4),
    @('     -', 'US Core Coverage Profile', '6.0.0', '', '', '', 3, 4),
    @('     -', 'US Core Observation Occupation Profile', '6.0.0', '', '', '', 3, 3),
    @('     -', 'US Core Observation Pregnancy Status Profile', '6.0.0', '', '', '', 3, 4),
    @('     -', 'US Core Observation Pregnancy Intent Profile', '6.0.0', '', '', '', 3, 4),
    @('     -', 'US Core Observation Clinical Result Profile', '6.0.0', '', '', 'Replaces US Core Observation Clinical Test Result Profile and US Core Observation Imaging Result Profile', 3, 4),
    @('     -', 'US Core Observation Screening Assessment Profile', '6.0.0', '', '', 'Replaces US Core Observation Survey Profile and US Core Observation SDOH Assessment Profile', 3, 3),
    @('     -', 'US Core MedicationDispense Profile', '6.0.0', '', '', '', 3, 4),
    @('     -', 'US Core Simple Observation Profile', '6.0.0', '', '', 'Replaces US Core Observation Social History Profile', 3, 4),
    @('     -', 'US Core Specimen Profile', '6.0.0', '', '', '', 3, 3),
    @('     -', 'US Core Average Blood Pressure Profile', '7.0.0', '', '', '', 2, 2),
    @('     -', 'US Core Care Experience Preference Profile', '7.0.0', '', '', '', 2, 2),
    @('     -', 'US Core Treatment Intervention Preference Profile', '7.0.0', '', '', '', 2, 2),
    @('     -', 'US Core Observation ADI Documentation Profile', '8.0.0', '', '', '', 2, 2),
    @('     -', 'US Core ADI DocumentReference Profile', '8.0.0', '', '', '', 2, 2),
    @('     -', 'US Core FamilyMemberHistory Profile', '9.0.0', '', '', '', 2, 2),
    @('     -', 'US Core PMO ServiceRequest Profile', '9.0.0', '', '', '', 2, 2),
    @('**Extensions**', '', '', '', '', '', '', ''),
    @('Sex of patient assigned at birth', 'US Core Birth Sex Extension', '3.1.1*', '', '', 'After version 6.0.0 this extension is no longer a USCDI requirement.', 3, 5),
    @('Argonaut ethnicity Extension', 'US Core Ethnicity Extension', '3.1.1*', '', '', '', 3, 5),
    @('Argonaut Race Extension', 'US Core Race Extension', '3.1.1*', '', '', '', 3, 5),
    @('     -', 'US Core Direct email Extension', '3.1.1*', '', '', '', 3, 3),
    @('     -', 'US Core Extension Questionnaire URI', '5.0.0', '', '', '', 3, 3),
    @('     -', 'US Core Gender Identity Extension', '5.0.0', '8.0.0', '', 'After version 7.0.0 this extension is no longer a USCDI requirement and has been deprecated.', 3, 3),
    @('     -', 'US Core Tribal Affiliation Extension', '6.0.0', '', '', '', 3, 3),
    @('     -', 'US Core Jurisdiction Extension', '6.0.0', '', '', '', 3, 2),
    @('     -', 'US Core USCDI Requirement Extension', '6.0.0', '', '', 'This extension is only used on US Core Profile StructureDefinition elements', 3, 5),
    @('     -', 'US Core Medication Adherence Extension', '7.0.0', '', '', '', 2, 2),
    @('     -', 'US Core Authentication Time Extension', '8.0.0', '9.0.0', '', 'Superseded  in version 9.0.0 by the Cross-version Extension for R5.DocumentReference.attester for use in FHIR R4.', 2, 2),
    @('     -', 'US Core Interpreter Needed Extension', '8.0.0', '', '', '', 2, 2),
    @('     -', 'US Core Sex Extension', '6.1.0', '8.0.1', '', 'Superseded  in version 8.0.1 by the US Core Individual Sex Extension.', 3, 2),
    @('     -', 'US Core Individual Sex Extension', '8.0.1', '', '', '', 'NA', 5),
    @('     -', 'US Core FamilyMemberHistory Recorder Extension', '9.0.0', '', '', '', 'NA', 2),
    @('', '', '', '', '', '\* 3.1.1  *or prior* version of US Core', '', '')
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
